$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'Sembol',
    'GOOGL',
    'GOOG',
    'AAPL',
    'MSFT',
    'AMZN',
    'META',
    'TSLA',
    'JPM',
    'V',
    'ORCL',
    'BAC',
    'MA',
    'BABA',
    'AMD',
    'PLTR',
    'NFLX',
    'GS',
    'JPM-PD',
    'JPM-PC',
    'UNH',
    'MS',
    'BAC-PK',
    'BAC-PL',
    'CSCO',
    'HSBC',
    'C',
    'NVO',
    'SAP',
    'IBM',
    'WFC',
    'BAC-PB',
    'AXP',
    'TMO',
    'CRM',
    'INTC',
    'WFC-PY',
    'DIS',
    'WFC-PL',
    'ISRG',
    'ABT',
    'SAN',
    'SCHW',
    'BX',
    'APP',
    'ANET',
    'SHOP',
    'ACN',
    'BLK',
    'UBER',
    'DHR',
    'BKNG',
    'QCOM',
    'HDB',
    'SPGI',
    'INTU',
    'UBS',
    'PDD',
    'BBVA',
    'NOW',
    'COF',
    'BSX',
    'SONY',
    'MDT',
    'PANW',
    'IBKR',
    'ADBE',
    'VRTX',
    'PGR',
    'CRWD',
    'SNY',
    'MELI',
    'KKR',
    'SPOT',
    'ADP',
    'CEG',
    'CVNA',
    'LYG',
    'ICE',
    'GD',
    'SNPS',
    'HOOD',
    'NKE',
    'BP',
    'MCO',
    'BCS',
    'MRSH',
    'DASH',
    'NU',
    'CDNS',
    'B',
    'NTES',
    'ELV',
    'ORLY',
    'ING',
    'MS-PK',
    'EMR',
    'BK',
    'MS-PI',
    'CRH',
    'ABNB',
    'TDG',
    'DB',
    'MS-PF',
    'MS-PE',
    'DELL',
    'APO',
    'GM',
    'INFY',
    'SE',
    'AON',
    'SNOW',
    'RELX',
    'NWG',
    'MRVL',
    'MS-PA',
    'LHX',
    'AJG',
    'TFC',
    'MFC',
    'DUK-PA',
    'RKT',
    'CTA-PB',
    'ADSK',
    'COIN',
    'NDAQ',
    'IDXX',
    'TRI',
    'BIDU',
    'PYPL',
    'CMG',
    'CCJ',
    'RBLX',
    'MET',
    'EA',
    'SCHW-PD',
    'WDAY',
    'CBRE',
    'EW',
    'ARES',
    'AXON',
    'DDOG',
    'ALNY',
    'ROK',
    'AMP',
    'HEI',
    'MSTR',
    'KGC',
    'MSCI',
    'TTWO',
    'SPG-PJ',
    'DAL',
    'WCN',
    'ROP',
    'JD',
    'EBAY',
    'RKLB',
    'MET-PA',
    'EL',
    'CTSH',
    'TCOM',
    'LVS',
    'IQV',
    'PUK',
    'RDDT',
    'CPRT',
    'XYZ',
    'HEI-A',
    'ALC',
    'GRMN',
    'MLM',
    'ASTS',
    'A',
    'PRU',
    'PAYX',
    'MDLN',
    'RMD',
    'FICO',
    'VEEV',
    'GEHC',
    'FISV',
    'TEAM',
    'RYAAY',
    'CPNG',
    'STT',
    'CLS',
    'SYM',
    'APO-PA',
    'XYL',
    'UAL',
    'SATS',
    'ZS',
    'INSM',
    'NTRA',
    'RJF',
    'MDB',
    'EXPE',
    'BE',
    'ESLT',
    'PSA-PK',
    'ALL-PH',
    'SOFI',
    'ALL-PB',
    'FOXA',
    'HUM',
    'WTW',
    'FIS',
    'FOX',
    'VRSK',
    'QSR',
    'FLUT',
    'ROL',
    'BNTX',
    'MTD',
    'SYF',
    'DXCM',
    'LPLA',
    'ZM',
    'HPE',
    'TME',
    'NTRS',
    'STLA',
    'CSGP',
    'BAP',
    'ALAB',
    'PAAS',
    'WIT',
    'HBAN',
    'EC',
    'BRO',
    'CFG',
    'EXE',
    'EFX',
    'FSLR',
    'STE',
    'OMC',
    'AVB',
    'VLTO',
    'DLR-PK',
    'CINF',
    'FCNCA',
    'WSM',
    'BR',
    'LDOS',
    'PSTG',
    'SQM',
    'FLEX',
    'AXIA-PC',
    'ILMN',
    'VRSN',
    'KEY',
    'TPG',
    'TROW',
    'WAT',
    'CRDO',
    'NBIS',
    'LULU',
    'OWL',
    'CNC',
    'FUTU',
    'AFRM',
    'DLR-PJ',
    'TLK',
    'CYBR',
    'FWONK',
    'ALB',
    'PSLV',
    'FWONA',
    'CG',
    'RBA',
    'CIB',
    'EXPD',
    'UTHR',
    'GPN',
    'CHRW',
    'BEKE',
    'SSNC',
    'GMAB',
    'PFG',
    'TWLO',
    'Q',
    'INCY',
    'HL',
    'CHKP',
    'PTC',
    'LTM',
    'NTAP',
    'TOST',
    'JBHT',
    'SGI',
    'GIB',
    'RIVN',
    'PODD',
    'TYL',
    'RVMD',
    'BWXT',
    'KTOS',
    'MRNA',
    'HIG-PG',
    'GRAB',
    'IOT',
    'FTV',
    'CX',
    'BURL',
    'U',
    'SMCI',
    'DKS',
    'HPQ',
    'STT-PG',
    'CRCL',
    'FITBI',
    'XPEV',
    'IT',
    'PSNYW',
    'NVT',
    'ALLY',
    'PNR',
    'PINS',
    'ESS',
    'WST',
    'HUBS',
    'NWS',
    'IREN',
    'ZG',
    'FN',
    'TRMB',
    'MEDP',
    'JLL',
    'Z',
    'TKO',
    'CDW',
    'RTO',
    'TRU',
    'TTD',
    'CDE',
    'NLY',
    'J',
    'HII',
    'MAA',
    'OKTA',
    'COO',
    'GFL',
    'WMG',
    'KSPI',
    'ROKU',
    'GEN',
    'NWSA',
    'EWBC',
    'FFIV',
    'IONQ',
    'ONON',
    'CRBG',
    'DKNG',
    'AA',
    'AVAV',
    'KEY-PI',
    'GH',
    'MLI',
    'HMY',
    'ULS',
    'GDDY',
    'ERIE',
    'ARCC',
    'ICLR',
    'FIG',
    'PNFP',
    'KRMN',
    'W',
    'TSEM',
    'CACI',
    'PEN',
    'EVR',
    'CELH',
    'DPZ',
    'BBY',
    'EMBJ',
    'SBSW',
    'LOGI',
    'EQH',
    'GWRE',
    'BILI',
    'RBRK',
    'HRL',
    'FIGR',
    'NBIX',
    'RVTY',
    'JKHY',
    'RYAN',
    'SF',
    'IONS',
    'PSKY',
    'CHWY',
    'UNM',
    'SNAP',
    'JEF',
    'OKLO',
    'HLI',
    'IVZ',
    'ACM',
    'BMNR',
    'GLXY',
    'AG',
    'GMED',
    'SNX',
    'DT',
    'ACGLO',
    'TXRH',
    'DOC',
    'SMMT',
    'RMBS',
    'EPAM',
    'JOBY',
    'CMA',
    'TEM',
    'FHN',
    'WYNN',
    'AEG',
    'AIZ',
    'EXEL',
    'NTNX',
    'BSY',
    'DOCU',
    'CPT',
    'BXP',
    'STN',
    'CEF',
    'MDGL',
    'BAH',
    'SCI',
    'MICC',
    'MP',
    'GL',
    'QGEN',
    'SARO',
    'CRL',
    'AR',
    'UHAL',
    'MOH',
    'DRS',
    'GME',
    'MANH',
    'FDS',
    'AFG',
    'CART',
    'SEIC',
    'TECH',
    'CAE',
    'YMM',
    'PCOR',
    'KLAR',
    'ARMK',
    'BAX',
    'PPC',
    'SSB',
    'CHYM',
    'VNO-PL',
    'UHAL-B',
    'SAIL',
    'SANM',
    'NGD',
    'GAP',
    'COMP',
    'WTFC',
    'BIO-B',
    'TTMI',
    'SOLS',
    'ARE',
    'BROS',
    'REXR',
    'AAL',
    'TTEK',
    'AYI',
    'UWMC',
    'GTLS',
    'ARWR',
    'EDU',
    'ONB',
    'SAIA',
    'RGEN',
    'DOX',
    'KNX',
    'STEP',
    'TFII',
    'MORN',
    'UEC',
    'LUMN',
    'GGAL',
    'QBTS',
    'TTAN',
    'RZB',
    'AGNCM',
    'AGNCN',
    'ZION',
    'PEGA',
    'VIPS',
    'AUR',
    'UGI',
    'AMTM',
    'PL'
)

for ($i = 0; $i -lt 499; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
